$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 18 (pushes old row18.. down by one) ---
$ws.Rows("18:18").Insert()

# Copy formatting for the new row18 (B:G) from row14, which has the matching
# visual pattern (fill/border/good-style) that the new row needs.
$ws.Range("B14:G14").Copy()
$ws.Range("B18:G18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Update core timer values ---
$ws.Range("D10").Value = 119
$ws.Range("D12").Value = 76

# PWM freq formula now divides by (D10+1) instead of D10 (preload enabled)
$ws.Range("D14").Formula = "=D5/(D8+1)/(D10+1)"

# --- PWM duty row (17) gains compare-event-duty columns F/G ---
$ws.Range("F17").Style = "Good"
$ws.Range("F17").Formula = "=F14*D17/100"
$ws.Range("G17").Value = "ms"

# --- New row 18: compare event period expressed in us ---
$ws.Range("F18").Formula = "=F17*1000"
$ws.Range("G18").Value = "us"

# --- Legend / measurement notes below the chart (rows 22-24) ---
$ws.Range("C22").Value = 119
$ws.Range("D22").Value = "'->"
$ws.Range("E22").Value = "1.25 us per"

$ws.Range("C23").Value = 38
$ws.Range("E23").Value = ".4us pos width"

$ws.Range("C24").Value = 76
$ws.Range("E24").Value = 0.798

# --- Column D needs an explicit width now that it holds more data ---
$ws.Columns("D").ColumnWidth = 10.59

# --- Selection / cursor position as last left by the author ---
$ws.Range("E24").Select()
